# tutorial covid model revised
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("simulator_scenarios")

# Update scenario values in row 2
$ws.Range("D2").Value = 5
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0

# Reset scrolled view / selection to cell F7
$ws.Range("F7").Select()

$wb.Save()
